$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 35945
$ws.Range("J109").Value = 35945
$ws.Range("L109").Value = 35945
$ws.Range("N109").Value = -38719

$ws.Range("H114").Value = 39602.5
$ws.Range("J114").Value = 39602.5
$ws.Range("L114").Value = 39602.5
$ws.Range("N114").Value = -48280.5

$ws.Range("H117").Value = 47031.5
$ws.Range("J117").Value = 47031.5
$ws.Range("L117").Value = 47031.5
$ws.Range("N117").Value = -56209.5

$ws.Range("H128").Value = 43221.4
$ws.Range("J128").Value = 43221.4
$ws.Range("L128").Value = 43221.4
$ws.Range("N128").Value = -53181.4

$ws.Range("H130").Value = 45450.668
$ws.Range("J130").Value = 45450.668
$ws.Range("L130").Value = 45450.668
$ws.Range("N130").Value = -55490.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H111").Value = 49632
$ws.Range("J111").Value = 49632
$ws.Range("L111").Value = 49632
$ws.Range("N111").Value = -57812

$ws.Range("H117").Value = 44364.285
$ws.Range("J117").Value = 44364.285
$ws.Range("L117").Value = 44364.285
$ws.Range("N117").Value = -53542.285

$ws.Range("H118").Value = 49409
$ws.Range("J118").Value = 49409
$ws.Range("L118").Value = 49409
$ws.Range("N118").Value = -52723

$ws.Range("H119").Value = 51500
$ws.Range("J119").Value = 51500
$ws.Range("L119").Value = 51500
$ws.Range("N119").Value = -61176

$ws.Range("H121").Value = 45255
$ws.Range("J121").Value = 45255
$ws.Range("L121").Value = 45255
$ws.Range("N121").Value = -48749

$ws.Range("H125").Value = 50694
$ws.Range("J125").Value = 50694
$ws.Range("L125").Value = 50694
$ws.Range("N125").Value = -60534

$ws.Range("H130").Value = 46210.5
$ws.Range("J130").Value = 46210.5
$ws.Range("L130").Value = 46210.5
$ws.Range("N130").Value = -56250.5

$ws.Range("H131").Value = 50277
$ws.Range("J131").Value = 50277
$ws.Range("L131").Value = 50277
$ws.Range("N131").Value = -60357

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 41072
$ws.Range("J108").Value = 41072
$ws.Range("L108").Value = 41072
$ws.Range("N108").Value = -48752

$ws.Range("H110").Value = 47383.332
$ws.Range("J110").Value = 47383.332
$ws.Range("L110").Value = 47383.332
$ws.Range("N110").Value = -55563.332

$ws.Range("H111").Value = 47694
$ws.Range("J111").Value = 47694
$ws.Range("L111").Value = 47694
$ws.Range("N111").Value = -55874

$ws.Range("H112").Value = 46484.332
$ws.Range("J112").Value = 46484.332
$ws.Range("L112").Value = 46484.332
$ws.Range("N112").Value = -49438.332

$ws.Range("H116").Value = 43448
$ws.Range("J116").Value = 43448
$ws.Range("L116").Value = 43448
$ws.Range("N116").Value = -52626

$ws.Range("H117").Value = 47387.75
$ws.Range("J117").Value = 47387.75
$ws.Range("L117").Value = 47387.75
$ws.Range("N117").Value = -56565.75

$ws.Range("H124").Value = 49992
$ws.Range("J124").Value = 49992
$ws.Range("L124").Value = 49992
$ws.Range("N124").Value = -59812

$ws.Range("H125").Value = 50772
$ws.Range("J125").Value = 50772
$ws.Range("L125").Value = 50772
$ws.Range("N125").Value = -60612

$ws.Range("H130").Value = 45085
$ws.Range("J130").Value = 45085
$ws.Range("L130").Value = 45085
$ws.Range("N130").Value = -55125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 42876.4
$ws.Range("J20").Value = 42876.4
$ws.Range("L20").Value = 42876.4
$ws.Range("N20").Value = -43348.4

$ws.Range("H30").Value = 42876.4
$ws.Range("J30").Value = 42876.4
$ws.Range("L30").Value = 42876.4
$ws.Range("N30").Value = -43058.4

$ws.Range("H111").Value = 48747
$ws.Range("J111").Value = 48747
$ws.Range("L111").Value = 48747
$ws.Range("N111").Value = -56927

$ws.Range("H112").Value = 28647.834
$ws.Range("J112").Value = 28647.834
$ws.Range("L112").Value = 28647.834
$ws.Range("N112").Value = -31601.834

$ws.Range("H116").Value = 36864.5
$ws.Range("J116").Value = 36864.5
$ws.Range("L116").Value = 36864.5
$ws.Range("N116").Value = -46042.5

$ws.Range("H119").Value = 46250.668
$ws.Range("J119").Value = 46250.668
$ws.Range("L119").Value = 46250.668
$ws.Range("N119").Value = -55926.668

$ws.Range("H128").Value = 42876.4
$ws.Range("J128").Value = 42876.4
$ws.Range("L128").Value = 42876.4
$ws.Range("N128").Value = -52836.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H110").Value = 47694
$ws.Range("J110").Value = 47694
$ws.Range("L110").Value = 47694
$ws.Range("N110").Value = -55874

$ws.Range("H114").Value = 43810.332
$ws.Range("J114").Value = 43810.332
$ws.Range("L114").Value = 43810.332
$ws.Range("N114").Value = -52488.332

$ws.Range("H116").Value = 39000
$ws.Range("J116").Value = 39000
$ws.Range("L116").Value = 39000
$ws.Range("N116").Value = -48178

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H130").Value = 45782.855
$ws.Range("J130").Value = 45782.855
$ws.Range("L130").Value = 45782.855
$ws.Range("N130").Value = -55822.855

$ws.Range("H139").Value = 32812.285
$ws.Range("J139").Value = 32812.285
$ws.Range("L139").Value = 32812.285
$ws.Range("N139").Value = -43092.285

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H108").Value = 48618
$ws.Range("J108").Value = 48618
$ws.Range("L108").Value = 48618
$ws.Range("N108").Value = -56298

$ws.Range("H110").Value = 40995
$ws.Range("J110").Value = 40995
$ws.Range("L110").Value = 40995
$ws.Range("N110").Value = -49175

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H114").Value = 39390
$ws.Range("J114").Value = 39390
$ws.Range("L114").Value = 39390
$ws.Range("N114").Value = -48068

$ws.Range("H116").Value = 47672
$ws.Range("J116").Value = 47672
$ws.Range("L116").Value = 47672
$ws.Range("N116").Value = -56850

$ws.Range("H119").Value = 47412
$ws.Range("J119").Value = 47412
$ws.Range("L119").Value = 47412
$ws.Range("N119").Value = -57088

$ws.Range("H120").Value = 54495
$ws.Range("J120").Value = 54495
$ws.Range("L120").Value = 54495
$ws.Range("N120").Value = -64171

$ws.Range("H125").Value = 49715
$ws.Range("J125").Value = 49715
$ws.Range("L125").Value = 49715
$ws.Range("N125").Value = -59555

$ws.Range("H127").Value = 50577.168
$ws.Range("J127").Value = 50577.168
$ws.Range("L127").Value = 50577.168
$ws.Range("N127").Value = -60497.168

$ws.Range("H130").Value = 47996
$ws.Range("J130").Value = 47996
$ws.Range("L130").Value = 47996
$ws.Range("N130").Value = -58036

$ws.Range("H133").Value = 22102.643
$ws.Range("J133").Value = 22102.643
$ws.Range("L133").Value = 22102.643
$ws.Range("N133").Value = -27162.643

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 47644
$ws.Range("J110").Value = 47644
$ws.Range("L110").Value = 47644
$ws.Range("N110").Value = -55824

$ws.Range("H117").Value = 47192
$ws.Range("J117").Value = 47192
$ws.Range("L117").Value = 47192
$ws.Range("N117").Value = -56370

$ws.Range("H120").Value = 45420
$ws.Range("J120").Value = 45420
$ws.Range("L120").Value = 45420
$ws.Range("N120").Value = -55096

$ws.Range("H121").Value = 44420
$ws.Range("J121").Value = 44420
$ws.Range("L121").Value = 44420
$ws.Range("N121").Value = -47914

$ws.Range("H128").Value = 49707
$ws.Range("J128").Value = 49707
$ws.Range("L128").Value = 49707
$ws.Range("N128").Value = -59667

$ws.Range("H131").Value = 47351.8
$ws.Range("J131").Value = 47351.8
$ws.Range("L131").Value = 47351.8
$ws.Range("N131").Value = -57431.8

Write-Output "applied 52 row updates"